# Add 2022-Q1 data:
#  - the current last sheet ("总计") is repurposed (renamed + refilled) as the
#    new "2022-Q1" per-fund detail sheet
#  - a fresh "总计" summary sheet is appended after it, cloned from the
#    "2021-Q4" sheet (for identical sheet-level formatting), cleared down and
#    refilled with the old summary rows plus a new leading row for 2022-Q1

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item(3)                 # "2021-Q4" detail sheet - style template
$q1New = $wb.Worksheets.Item($wb.Worksheets.Count)   # currently "总计" (sheetId 4)

# --- rename + rebuild the old "总计" sheet into the "2022-Q1" detail sheet ---
$q1New.Name = "2022-Q1"

# Pull the header + index-column formatting from the existing 2021-Q4 sheet so
# the new sheet matches the established per-quarter layout (bold/boxed style).
$q4.Range("B1:H1").Copy()
$q1New.Range("B1:H1").PasteSpecial(-4122)

$q4.Range("A2:A4").Copy()
$q1New.Range("A2:A4").PasteSpecial(-4122)

$q1New.Range("B1").Value = "基金代码"
$q1New.Range("C1").Value = "基金名称"
$q1New.Range("D1").Value = "基金规模"
$q1New.Range("E1").Value = "股票总仓位"
$q1New.Range("F1").Value = "仓位占比"
$q1New.Range("G1").Value = "持有市值(亿元)"
$q1New.Range("H1").Value = "仓位排名"

$q1New.Range("A2").Value = 0
$q1New.Range("A3").Value = 1
$q1New.Range("A4").Value = 2

# Numeric-looking text fields must stay text (leading zeros in codes, and the
# source data stores these figures as text, not numbers).
function Set-TextValue($rng, $value) {
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue $q1New.Range("B2") "010490"
Set-TextValue $q1New.Range("B3") "009023"
Set-TextValue $q1New.Range("B4") "010491"

$q1New.Range("C2").Value = "鹏华高质量增长混合A"
$q1New.Range("C3").Value = "鹏华稳健回报混合"
$q1New.Range("C4").Value = "鹏华高质量增长混合C"

Set-TextValue $q1New.Range("D2") "13.31"
Set-TextValue $q1New.Range("E2") "93.61"
Set-TextValue $q1New.Range("F2") "3.89"
Set-TextValue $q1New.Range("G2") "0.5178"

Set-TextValue $q1New.Range("D3") "3.52"
Set-TextValue $q1New.Range("E3") "93.91"
Set-TextValue $q1New.Range("F3") "6.14"
Set-TextValue $q1New.Range("G3") "0.2161"

Set-TextValue $q1New.Range("D4") "0.28"
Set-TextValue $q1New.Range("E4") "93.61"
Set-TextValue $q1New.Range("F4") "3.89"
Set-TextValue $q1New.Range("G4") "0.0109"

$q1New.Range("H2").Value = 5
$q1New.Range("H3").Value = 4
$q1New.Range("H4").Value = 5

# --- append a brand-new "总计" summary sheet after "2022-Q1" ---
# Clone the "2021-Q4" sheet so the new sheet inherits identical sheet-level
# properties (sheetPr / sheetFormatPr / pageMargins), then trim it back down
# to the small 4-column summary shape and refill its values.
$q4.Copy($null, $q1New)
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Name = "总计"

$total.Range("E1:H13").Clear()
$total.Range("A6:D13").Clear()

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.74

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 12
$total.Range("D3").Value = 2.76

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 6
$total.Range("D4").Value = 1.3

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.02

# Restore selection to the first sheet (matches the original file's state,
# and drops the transient "active sheet" marker the copy left behind).
$wb.Worksheets.Item(1).Activate()
